# Rename speaker labels in column D:
#  - "HILLARY LEWIS-WOLFSEN" -> "T"
#  - "STUDENT A" / "STUDENT B" -> "S"
# Other speaker names (e.g. KEVIN, ANDREW, SI) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()

    if ($val -eq "HILLARY LEWIS-WOLFSEN") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT A" -or $val -eq "STUDENT B") {
        $cell.Value = "S"
    }
}
